# Regenerate the "K" column (column G) values in the save-data sheet.
# The commit replaces the old Strike#-derived K values with newly
# calculated s_vals for each row (rows 2-34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 4
    4  = 1
    5  = 2
    6  = 5
    7  = 2
    8  = 5
    9  = 0
    10 = 5
    11 = 2
    12 = 5
    13 = 4
    14 = 2
    15 = 2
    16 = 4
    17 = 2
    18 = 2
    19 = 4
    20 = 4
    21 = 4
    22 = 3
    23 = 2
    24 = 4
    25 = 6
    26 = 1
    27 = 3
    28 = 1
    29 = 5
    30 = 3
    31 = 3
    32 = 4
    33 = 3
    34 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
